$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.367430090904236
$ws.Range("B1").Value = 1.599879860877991
$ws.Range("C1").Value = 1.358075857162476
$ws.Range("D1").Value = 1.411431789398193
$ws.Range("E1").Value = 0.9763326048851013
